$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.407.41"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.838.24"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6244"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07391"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2936"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07671"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.849.73"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.017"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6743"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009165"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.888"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "29.383.80"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "2.093.73"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.386"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1406"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.489"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.233"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.117"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.097"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.142"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7236"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.617"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.895"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.72%  "
$ws.Range("D39").Value = "1.222.72"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01762"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9137"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "2.009.18"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5073"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.227"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4050"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1151"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.98%  "
